$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 17888.666
$ws.Range("I21").Value = 9000
$ws.Range("J21").Value = 18999.75
$ws.Range("K21").Value = 9000
$ws.Range("L21").Value = 18999.75
$ws.Range("M21").Value = -8532
$ws.Range("N21").Value = -19935.75

$ws.Range("H23").Value = 17888.666
$ws.Range("I23").Value = 9000
$ws.Range("J23").Value = 18999.75
$ws.Range("K23").Value = 9000
$ws.Range("L23").Value = 18999.75
$ws.Range("M23").Value = -8766
$ws.Range("N23").Value = -19467.75

$ws.Range("H29").Value = 30
$ws.Range("I29").Value = 30
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 90
$ws.Range("L29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = 191

$ws.Range("H32").Value = 692.75
$ws.Range("I32").Value = 499.66666
$ws.Range("J32").Value = 808.6
$ws.Range("K32").Value = 499.66666
$ws.Range("L32").Value = 808.6
$ws.Range("M32").Value = -173.66666
$ws.Range("N32").Value = -1460.6

$ws.Range("H38").Value = 668.2
$ws.Range("I38").Value = 235.47058
$ws.Range("J38").Value = 1587.75
$ws.Range("K38").Value = 706.41174
$ws.Range("L38").Value = 4763.25
$ws.Range("M38").Value = -334.41174
$ws.Range("N38").Value = -5507.25

$ws.Range("H58").Value = 1205.5555
$ws.Range("I58").Value = 808.3333
$ws.Range("J58").Value = 2000
$ws.Range("K58").Value = 2424.9999
$ws.Range("L58").Value = 6000
$ws.Range("M58").Value = -2274.9999
$ws.Range("N58").Value = -6300

$ws.Range("H87").Value = 25375.691
$ws.Range("J87").Value = 25375.691
$ws.Range("L87").Value = 25375.691
$ws.Range("N87").Value = -27871.691

$ws.Range("H90").Value = 25375.691
$ws.Range("J90").Value = 25375.691
$ws.Range("L90").Value = 76127.073
$ws.Range("N90").Value = -88607.073

$ws.Range("H98").Value = 1287.7778
$ws.Range("I98").Value = 1136.25
$ws.Range("K98").Value = 1136.25
$ws.Range("M98").Value = 361.75

$ws.Range("H122").Value = 1287.7778
$ws.Range("I122").Value = 1136.25
$ws.Range("K122").Value = 3408.75
$ws.Range("M122").Value = -958.75

$ws.Range("H124").Value = 80990
$ws.Range("J124").Value = 80990
$ws.Range("L124").Value = 80990
$ws.Range("N124").Value = -90810

$ws.Range("H125").Value = 9617180
$ws.Range("I125").Value = 41667430
$ws.Range("J125").Value = 2104.4
$ws.Range("K125").Value = 375006870
$ws.Range("L125").Value = 18939.6
$ws.Range("M125").Value = -375004410
$ws.Range("N125").Value = -23859.6

$ws.Range("H132").Value = 5282.3687
$ws.Range("I132").Value = 5398.5
$ws.Range("J132").Value = 4663
$ws.Range("K132").Value = 16195.5
$ws.Range("L132").Value = 13989
$ws.Range("M132").Value = -13665.5
$ws.Range("N132").Value = -19049

$ws.Range("H138").Value = 3461.73
$ws.Range("I138").Value = 1390.75
$ws.Range("J138").Value = 3856.2024
$ws.Range("K138").Value = 4172.25
$ws.Range("L138").Value = 11568.6072
$ws.Range("M138").Value = 967.75
$ws.Range("N138").Value = -21848.6072

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").ClearContents()
$ws.Range("N52").Value = 0

$ws.Range("H74").Value = 13891776
$ws.Range("I74").Value = 2201
$ws.Range("J74").Value = 35718252
$ws.Range("K74").Value = 2201
$ws.Range("L74").Value = 35718252
$ws.Range("M74").Value = -1327
$ws.Range("N74").Value = -35720000

$ws.Range("H77").Value = 13891776
$ws.Range("I77").Value = 2201
$ws.Range("J77").Value = 35718252
$ws.Range("K77").Value = 11005
$ws.Range("L77").Value = 178591260
$ws.Range("M77").Value = -6637
$ws.Range("N77").Value = -178599996

$ws.Range("H135").Value = 79619.336
$ws.Range("J135").Value = 79619.336
$ws.Range("L135").Value = 79619.336
$ws.Range("N135").Value = -89759.336

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 33585
$ws.Range("J52").Value = 33585
$ws.Range("L52").Value = 33585
$ws.Range("N52").Value = -34173

$ws.Range("H132").Value = 23811998
$ws.Range("I132").Value = 25002214
$ws.Range("K132").Value = 75006642
$ws.Range("M132").Value = -75004112

$ws.Range("H134").Value = 10006733
$ws.Range("I134").Value = 13165038
$ws.Range("K134").Value = 39495114
$ws.Range("M134").Value = -39492579

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 901.431
$ws.Range("I5").Value = 677.9722
$ws.Range("J5").Value = 1267.091
$ws.Range("K5").Value = 2033.9166
$ws.Range("L5").Value = 3801.273
$ws.Range("M5").Value = -1921.9166
$ws.Range("N5").Value = -4025.273

$ws.Range("H12").Value = 63.434784
$ws.Range("I12").Value = 29.666666
$ws.Range("J12").Value = 126.75
$ws.Range("K12").Value = 88.99999800000001
$ws.Range("L12").Value = 380.25
$ws.Range("M12").Value = 84.00000199999999
$ws.Range("N12").Value = -726.25

$ws.Range("H132").Value = 2736.3845
$ws.Range("I132").Value = 2447.25
$ws.Range("K132").Value = 22025.25
$ws.Range("M132").Value = -19495.25

$ws.Range("H135").Value = 901.431
$ws.Range("I135").Value = 677.9722
$ws.Range("J135").Value = 1267.091
$ws.Range("K135").Value = 6101.749800000001
$ws.Range("L135").Value = 11403.819
$ws.Range("M135").Value = -3566.749800000001
$ws.Range("N135").Value = -16473.819

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").ClearContents()
$ws.Range("N42").Value = 0

$ws.Range("H46").Value = 4923
$ws.Range("J46").Value = 5249.9165
$ws.Range("L46").Value = 5249.9165
$ws.Range("N46").Value = -5561.9165

$ws.Range("H80").Value = 2351475
$ws.Range("I80").Value = 4501450
$ws.Range("J80").Value = 201500
$ws.Range("K80").Value = 4501450
$ws.Range("L80").Value = 201500
$ws.Range("M80").Value = -4500452
$ws.Range("N80").Value = -203496

$ws.Range("H82").Value = 31999
$ws.Range("J82").Value = 31999
$ws.Range("L82").Value = 31999
$ws.Range("N82").Value = -32765

$ws.Range("H83").Value = 2351475
$ws.Range("I83").Value = 4501450
$ws.Range("J83").Value = 201500
$ws.Range("K83").Value = 22507250
$ws.Range("L83").Value = 1007500
$ws.Range("M83").Value = -22502258
$ws.Range("N83").Value = -1017484

$ws.Range("H85").Value = 31999
$ws.Range("J85").Value = 31999
$ws.Range("L85").Value = 31999
$ws.Range("N85").Value = -34651

$ws.Range("H94").Value = 25448
$ws.Range("J94").Value = 25448
$ws.Range("L94").Value = 25448
$ws.Range("N94").Value = -26800

$ws.Range("H107").Value = 420.69232
$ws.Range("I107").Value = 350
$ws.Range("J107").Value = 579.75
$ws.Range("K107").Value = 350
$ws.Range("L107").Value = 579.75
$ws.Range("M107").Value = 1570
$ws.Range("N107").Value = -4419.75

$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").ClearContents()
$ws.Range("N115").Value = 0

$ws.Range("H121").Value = 59500
$ws.Range("J121").Value = 59500
$ws.Range("L121").Value = 59500
$ws.Range("N121").Value = -62994

$ws.Range("H132").Value = 38468056
$ws.Range("I132").Value = 52638656
$ws.Range("J132").Value = 5001.5713
$ws.Range("K132").Value = 157915968
$ws.Range("L132").Value = 15004.7139
$ws.Range("M132").Value = -157913438
$ws.Range("N132").Value = -20064.7139

$ws.Range("H137").Value = 72797.14
$ws.Range("J137").Value = 72797.14
$ws.Range("L137").Value = 72797.14
$ws.Range("N137").Value = -82997.14

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7818.933
$ws.Range("I7").Value = 7208.4
$ws.Range("K7").Value = 7208.4
$ws.Range("M7").Value = -7096.4

$ws.Range("H126").Value = 7818.933
$ws.Range("I126").Value = 7208.4
$ws.Range("K126").Value = 21625.2
$ws.Range("M126").Value = -19155.2
